# Using Kubernetes Day 5.pptx - edits:
#  1. Refresh the cached "datetimeFigureOut" date placeholder text
#     (1/2/22 -> 3/3/22) on the Slide Master and every Slide Layout.
#  2. Fix the title textbox on slide 1 from "Day 6" to "Day 5".

$p = $ppt.ActivePresentation

$oldDate = "1/2/22"
$newDate = "3/3/22"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if (-not $shp.HasTextFrame) { continue }
        $isDate = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDate = $true
            }
        } catch {
            $isDate = $false
        }
        if (-not $isDate) { continue }

        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq $oldDate) {
            $tr.Text = $newDate
        }
    }
}

# Slide Master date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every Slide Layout's date placeholder.
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}

# Slide 1: "Day 6" -> "Day 5" textbox (title card).
# The box holds two runs ("Day " + "6"); drop the first run entirely so
# the remaining text no longer shares a prefix with the replacement,
# then retype it fully -- this collapses the paragraph into a single
# run (matching a full retype of the whole textbox).
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shp = $slide1.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "Day 6") {
        $tr = $shp.TextFrame.TextRange
        $tr.Characters(1, 4).Text = ""
        $tr.Text = "Day 5"
    }
}
